# Auto-generated edit script: updates crypto price/volume table
# to match the 'Updated cryptos list' GitHub Actions commit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '70.883.92'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +7.07%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.635.03'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +6.75%  '
$ws.Range("E4").Value = '  +0.04%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '593.51'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +4.88%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '192.87'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +8.55%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.654'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +3.48%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '3.623.38'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +6.48%  '
$ws.Range("E9").Value = '  +0.11%  '
$ws.Range("E10").Value = '  +2.25%  '
$ws.Range("E11").Value = '  +4.55%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '58.21'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +7.17%  '
$ws.Range("E13").Value = '  +6.27%  '
$ws.Range("E14").Value = '  +5.44%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '4.217.63'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +7.20%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '3.634.23'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +6.66%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '19.42'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +6.22%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '70.794.74'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +7.18%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '12.62'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +5.57%  '
$ws.Range("E20").Value = '  +0.72%  '
$ws.Range("E21").Value = '  +5.21%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '495.69'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +6.59%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '5.43'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +9.42%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '17.12'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +14.95%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '4.51'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +9.15%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '91.21'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +1.30%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '3.12'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +6.40%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '11.29'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +5.25%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '9.47'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +7.78%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '32.45'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +3.64%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '7.61'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +14.01%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '12.25'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +6.26%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '619.98'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +6.54%  '
$ws.Range("E34").Value = '  +8.45%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '65.33'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +4.37%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.0₃0833'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +9.80%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.415'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +9.07%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.148'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +3.73%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '38.27'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +5.38%  '
$ws.Range("E40").Value = '  +0.00%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '3.68'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +2.94%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '3.334.12'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +6.62%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '3.08'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +7.43%  '
$ws.Range("E44").Value = '  +7.04%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '2.71'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +9.34%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '3.34'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +4.72%  '
$ws.Range("E47").Value = '  +3.07%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '9.27'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +8.47%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '2.74'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +6.07%  '
$ws.Range("E50").Value = '  +4.82%  '
$ws.Range("B51").Value = 'Monero'
$ws.Range("C51").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '143.49'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +1.63%  '
